$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9423993556369226
$ws.Range("C2").Value = 0.0698908199847601
$ws.Range("D2").Value = 0.07693026849750595
$ws.Range("E2").Value = 0.09638725095160083
$ws.Range("G2").Value = 2.813795971865773
$ws.Range("H2").Value = 2.144975901200553
$ws.Range("K2").Value = 0.5032953942100278
$ws.Range("L2").Value = 0.2118762749613836
$ws.Range("B3").Value = 0.9140284257348412
$ws.Range("C3").Value = 0.06812135611863823
$ws.Range("D3").Value = 0.07001362871633887
$ws.Range("E3").Value = 0.09528090821886082
$ws.Range("G3").Value = 2.737642155308947
$ws.Range("H3").Value = 2.112649597458642
$ws.Range("K3").Value = 0.4780125753428024
$ws.Range("L3").Value = 0.205756322730636
$ws.Range("B4").Value = 0.8973565204716749
$ws.Range("C4").Value = 0.06700015309624163
$ws.Range("D4").Value = 0.06580612646935435
$ws.Range("E4").Value = 0.094638173090555
$ws.Range("G4").Value = 2.691600719588251
$ws.Range("H4").Value = 2.0933036386391
$ws.Range("K4").Value = 0.4629546324703
$ws.Range("L4").Value = 0.202132386448568
$ws.Range("B5").Value = 0.8907504930698167
$ws.Range("C5").Value = 0.06653444758534732
$ws.Range("D5").Value = 0.06410133258202677
$ws.Range("E5").Value = 0.09438544283780104
$ws.Range("G5").Value = 2.6730176061007
$ws.Range("H5").Value = 2.085545894753864
$ws.Range("K5").Value = 0.4569352466074719
$ws.Range("L5").Value = 0.2006891732766007
$ws.Range("B6").Value = 0.8896649117655215
$ws.Range("C6").Value = 0.06645658364534412
$ws.Range("D6").Value = 0.06381884182827946
$ws.Range("E6").Value = 0.09434403206698505
$ws.Range("G6").Value = 2.669942678449956
$ws.Range("H6").Value = 2.084265317124505
$ws.Range("K6").Value = 0.4559427814169652
$ws.Range("L6").Value = 0.200451554694709
$ws.Range("B7").Value = 0.8972666686598245
$ws.Range("C7").Value = 0.06699390817516715
$ws.Range("D7").Value = 0.06578309547759886
$ws.Range("E7").Value = 0.09463472747956558
$ws.Range("G7").Value = 2.691349377499506
$ws.Range("H7").Value = 2.093198505905349
$ws.Range("K7").Value = 0.4628729801839881
$ws.Range("L7").Value = 0.2021127869221999
$ws.Range("B8").Value = 0.9324617012824774
$ws.Range("C8").Value = 0.06928787821234295
$ws.Range("D8").Value = 0.07453718562732092
$ws.Range("E8").Value = 0.09599819171637236
$ws.Range("G8").Value = 2.787388490515269
$ws.Range("H8").Value = 2.133725196366953
$ws.Range("K8").Value = 0.4944810165004583
$ws.Range("L8").Value = 0.2097383191214561
$ws.Range("B9").Value = 1.007427353539441
$ws.Range("C9").Value = 0.07351387577449486
$ws.Range("D9").Value = 0.09202134559599529
$ws.Range("E9").Value = 0.09896251183224081
$ws.Range("G9").Value = 2.981478650703508
$ws.Range("H9").Value = 2.217211627101619
$ws.Range("K9").Value = 0.5601776239616925
$ws.Range("L9").Value = 0.2257569478733501
$ws.Range("B10").Value = 1.066158091266431
$ws.Range("C10").Value = 0.07645719414164631
$ws.Range("D10").Value = 0.1050690870201834
$ws.Range("E10").Value = 0.1013185132157588
$ws.Range("G10").Value = 3.127693844187604
$ws.Range("H10").Value = 2.281039972787795
$ws.Range("K10").Value = 0.6107392781119074
$ws.Range("L10").Value = 0.2381820158907004
$ws.Range("B11").Value = 1.093676111492044
$ws.Range("C11").Value = 0.07776210814053286
$ws.Range("D11").Value = 0.1110506846526675
$ws.Range("E11").Value = 0.1024292319185989
$ws.Range("G11").Value = 3.195020427563463
$ws.Range("H11").Value = 2.310627992625655
$ws.Range("K11").Value = 0.6342464188481358
$ws.Range("L11").Value = 0.2439786078550128
$ws.Range("B12").Value = 1.104212055211065
$ws.Range("C12").Value = 0.07825143083416464
$ws.Range("D12").Value = 0.1133225214715452
$ws.Range("E12").Value = 0.1028554476097128
$ws.Range("G12").Value = 3.22063372076542
$ws.Range("H12").Value = 2.321912267375012
$ws.Range("K12").Value = 0.6432212420875771
$ws.Range("L12").Value = 0.2461944873331987
$ws.Range("B13").Value = 1.101937811413961
$ws.Range("C13").Value = 0.07814625966095434
$ws.Range("D13").Value = 0.1128329395213541
$ws.Range("E13").Value = 0.1027634047704993
$ws.Range("G13").Value = 3.21511216250542
$ws.Range("H13").Value = 2.319478433323411
$ws.Range("K13").Value = 0.641285092803372
$ws.Range("L13").Value = 0.2457163301240257
$ws.Range("B14").Value = 1.09454059422211
$ws.Range("C14").Value = 0.07780246119965994
$ws.Range("D14").Value = 0.1112374544410528
$ws.Range("E14").Value = 0.1024641844479248
$ws.Range("G14").Value = 3.197125272495612
$ws.Range("H14").Value = 2.311554751654057
$ws.Range("K14").Value = 0.6349833142334376
$ws.Range("L14").Value = 0.2441604915669302
$ws.Range("B15").Value = 1.090024630423926
$ws.Range("C15").Value = 0.07759124915123294
$ws.Range("D15").Value = 0.1102610540557265
$ws.Range("E15").Value = 0.1022816342712538
$ws.Range("G15").Value = 3.186123212880318
$ws.Range("H15").Value = 2.306711693402178
$ws.Range("K15").Value = 0.6311328338763929
$ws.Range("L15").Value = 0.2432102112730519
$ws.Range("B16").Value = 1.064375866965662
$ws.Range("C16").Value = 0.0763712365602558
$ws.Range("D16").Value = 0.10467911238689
$ws.Range("E16").Value = 0.1012467098212397
$ws.Range("G16").Value = 3.123310358050702
$ws.Range("H16").Value = 2.279117482687354
$ws.Range("K16").Value = 0.6092132516726849
$ws.Range("L16").Value = 0.2378061064615622
$ws.Range("B17").Value = 1.048846552133199
$ws.Range("C17").Value = 0.07561413679336937
$ws.Range("D17").Value = 0.101266656929127
$ws.Range("E17").Value = 0.1006218018160432
$ws.Range("G17").Value = 3.084985766765925
$ws.Range("H17").Value = 2.2623311316475
$ws.Range("K17").Value = 0.5958962503275131
$ws.Range("L17").Value = 0.2345278924600507
$ws.Range("B18").Value = 1.039989888785612
$ws.Range("C18").Value = 0.07517546695838462
$ws.Range("D18").Value = 0.09930823701768077
$ws.Range("E18").Value = 0.1002660380387219
$ws.Range("G18").Value = 3.063018877133004
$ws.Range("H18").Value = 2.252728008794463
$ws.Range("K18").Value = 0.5882843050713973
$ws.Range("L18").Value = 0.2326559342260452
$ws.Range("B19").Value = 1.037004114724112
$ws.Range("C19").Value = 0.07502638806859352
$ws.Range("D19").Value = 0.09864589093093912
$ws.Range("E19").Value = 0.1001462120723247
$ws.Range("G19").Value = 3.055594347356788
$ws.Range("H19").Value = 2.249485461458221
$ws.Range("K19").Value = 0.5857152028726773
$ws.Range("L19").Value = 0.2320244520646781
$ws.Range("B20").Value = 1.05049186996618
$ws.Range("C20").Value = 0.07569506265461001
$ws.Range("D20").Value = 0.101629469308179
$ws.Range("E20").Value = 0.1006879447852285
$ws.Range("G20").Value = 3.089057568046513
$ws.Range("H20").Value = 2.264112687753709
$ws.Range("K20").Value = 0.5973089345423546
$ws.Range("L20").Value = 0.23487545752225
$ws.Range("B21").Value = 1.096710200296371
$ws.Range("C21").Value = 0.07790357338139842
$ws.Range("D21").Value = 0.1117059034436068
$ws.Range("E21").Value = 0.1025519203741148
$ws.Range("G21").Value = 3.202405240135988
$ws.Range("H21").Value = 2.313879955399898
$ws.Range("K21").Value = 0.6368323103679359
$ws.Range("L21").Value = 0.2446169129675582
$ws.Range("B22").Value = 1.127589712117697
$ws.Range("C22").Value = 0.0793189021037648
$ws.Range("D22").Value = 0.1183307452702138
$ws.Range("E22").Value = 0.1038028431162417
$ws.Range("G22").Value = 3.277173949173118
$ws.Range("H22").Value = 2.346871906931256
$ws.Range("K22").Value = 0.6630898404124821
$ws.Range("L22").Value = 0.2511050016161391
$ws.Range("B23").Value = 1.111047049578701
$ws.Range("C23").Value = 0.07856605834651731
$ws.Range("D23").Value = 0.1147913108684691
$ws.Range("E23").Value = 0.103132206703382
$ws.Range("G23").Value = 3.237204954214008
$ws.Range("H23").Value = 2.329220650740353
$ws.Range("K23").Value = 0.6490365357529129
$ws.Range("L23").Value = 0.2476310450696388
$ws.Range("B24").Value = 1.049747799560606
$ws.Range("C24").Value = 0.07565848667958619
$ws.Range("D24").Value = 0.1014654311073144
$ws.Range("E24").Value = 0.1006580306350351
$ws.Range("G24").Value = 3.08721649987254
$ws.Range("H24").Value = 2.263307098077178
$ws.Range("K24").Value = 0.5966701224347446
$ws.Range("L24").Value = 0.234718283710734
$ws.Range("B25").Value = 0.9865075811472082
$ws.Range("C25").Value = 0.07239932861846654
$ws.Range("D25").Value = 0.08725654825930462
$ws.Range("E25").Value = 0.09812937740191074
$ws.Range("G25").Value = 2.928344677287981
$ws.Range("H25").Value = 2.194192109971027
$ws.Range("K25").Value = 0.5420042845127284
$ws.Range("L25").Value = 0.2213087602360133
